$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fix the visual layout of the "Nama Bank" column: indent + left align ---
# Apply directly to C2 first (a single cell clone keeps the style table tidy),
# then format-paint the same look onto the rest of the column / new cells.
$c2 = $ws.Cells.Item(2, 3)
$c2.IndentLevel = 1
$c2.HorizontalAlignment = -4131

$c2.Copy()
$ws.Range("C3:C9").PasteSpecial(-4122)
$c2.Copy()
$ws.Range("D6:D9").PasteSpecial(-4122)
$c2.Copy()
$ws.Range("E6:E9").PasteSpecial(-4122)

# --- Copy the look of the existing "Nama"/"NIK" columns down onto the new rows ---
$ws.Range("A2:B5").Copy()
$ws.Range("A6:B9").PasteSpecial(-4122)

# --- Add the new account rows (same users as rows 2-5, now with a real Nomor Rekening) ---
$rows = @(
    @{ A = "Administrator"; B = "20220001J"; C = "PANIN";    D = "Admin Administrator "; E = 982374982374 },
    @{ A = "Admin HRD";     B = "20221111J"; C = "PERMATA";  D = "Administrator HRD";    E = 6723942389 },
    @{ A = "Approver";      B = "20220002J"; C = "OCBC";     D = "Manager Approver";     E = 932402093 },
    @{ A = "User";          B = "20220003J"; C = "DANAMON";  D = "User Pegawai";         E = 8792347234 }
)

$r = 6
foreach ($row in $rows) {
    $ws.Cells.Item($r, 1).Value = $row.A
    $ws.Cells.Item($r, 2).Value = $row.B
    $ws.Cells.Item($r, 3).Value = $row.C
    $ws.Cells.Item($r, 4).Value = $row.D
    $ws.Cells.Item($r, 5).Value = $row.E
    $r++
}

$ws.Range("E10").Select()
